$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("J18").Value = "3.1/2"""
$ws.Range("J19").Value = "3.1/2"""
$ws.Range("J34").Value = "3.1/2"""
$ws.Range("J61").Value = "3000 PSI"
$ws.Range("L68").Value = "SPECIAL"
$ws.Range("J81").NumberFormat = "@"
$ws.Range("J81").Value = "27"
$ws.Range("J81").ClearFormats()
$ws.Range("J82").NumberFormat = "@"
$ws.Range("J82").Value = "27"
$ws.Range("J82").ClearFormats()
$ws.Range("L84").Value = "NTZ 400"
$ws.Range("L85").Value = "NTZ 400"
$ws.Range("L86").Value = "NETZSCH"
$ws.Range("J103").NumberFormat = "@"
$ws.Range("J103").Value = "27"
$ws.Range("J103").ClearFormats()
$ws.Range("K110").Value = "CI"
$ws.Range("K125").Value = "SS"
$ws.Range("K126").Value = "St"
$ws.Range("K127").Value = "St"
$ws.Range("K129").Value = "St"
$ws.Range("K130").Value = "St"
$ws.Range("K131").Value = "St"
$ws.Range("K132").Value = "St"
$ws.Range("K133").Value = "St"
$ws.Range("K136").Value = "St"
$ws.Range("J137").NumberFormat = "@"
$ws.Range("J137").Value = "90"
$ws.Range("J137").ClearFormats()
$ws.Range("K138").Value = "St"
$ws.Range("K139").Value = "St"
$ws.Range("K140").Value = "St"
$ws.Range("K141").Value = "St"
$ws.Range("K144").Value = "SS"
$ws.Range("K145").Value = "St"
$ws.Range("K146").Value = "St"
$ws.Range("K147").Value = "St"
$ws.Range("K151").Value = "SS"
$ws.Range("K153").Value = "PE"
$ws.Range("J154").NumberFormat = "@"
$ws.Range("J154").Value = "30"
$ws.Range("J154").ClearFormats()
$ws.Range("J155").NumberFormat = "@"
$ws.Range("J155").Value = "20"
$ws.Range("J155").ClearFormats()
$ws.Range("K156").Value = "St"
$ws.Range("K157").Value = "St"
$ws.Range("J161").Value = "1.1/4"""
$ws.Range("K161").Value = "SAE1045"
$ws.Range("L161").Value = "API 6A"
$ws.Range("J162").Value = "1.1/4"""
$ws.Range("K162").Value = "SAE1045"
$ws.Range("L162").Value = "API 6A"
$ws.Range("L163").NumberFormat = "@"
$ws.Range("L163").Value = "299"
$ws.Range("L163").ClearFormats()
$ws.Range("J166").Value = "1.1/4"""
$ws.Range("L166").Value = "API 6A"
$ws.Range("J177").NumberFormat = "@"
$ws.Range("J177").Value = "27"
$ws.Range("J177").ClearFormats()
$ws.Range("J181").Value = "1.1/4"""
$ws.Range("L181").Value = "API 6A"
$ws.Range("L182").Value = "NE"
$ws.Range("J184").Value = "4.1/2"""
$ws.Range("K184").Value = "GGG50"
